$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 51; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = (-$b)
    $ws.Cells.Item($r, 5).Value = ($b * $b)
}

$ws.Range("C52").Formula = "=SUM(D2:D51)"
$ws.Range("E52").Formula = "=SUM(E2:E51)"
$ws.Range("E53").Formula = "=AVERAGE(E2:E51)"
